$wb = $excel.ActiveWorkbook

# --- Update "Status" values from "Ready for handoff" to "In Translation" ---

# Overview sheet: zh-cn (col E) / de-de (col F) status columns, rows 2-3
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F3").Value = "In Translation"

# zh-cn sheet: Status column (col C), rows 2-3
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C3").Value = "In Translation"

# de-de sheet: Status column (col C), rows 2-3
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C3").Value = "In Translation"

# --- Narrow the "Status" columns to fit the shorter replacement text ---
# (originally sized for "Ready for handoff", now re-fit for "In Translation")
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
